$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new tasks inserted under the "Implementacija" block
$ws.Range("B30").Value = "Povezivanje frontend i backend dijela "
$ws.Range("B31").Value = "Povezivanje baze podataka s aplikacijom"

# "Testiranje" (previously on B33) moves down to B35, keeping its bold style
$ws.Range("B35").Value = "Testiranje"
$ws.Range("B35").Font.Bold = $true

# B33 keeps its (bold) style but loses its text
$ws.Range("B33").Value = ""

# New deploy task appended as a new row
$ws.Range("B37").Value = "Deploy sustava (stavljen ovdje a ne pod implementaciju jer kao necemo pustat u pogon dok se ne testira sve)"

# Column width tweaks
$ws.Columns.Item(2).ColumnWidth = 43
$ws.Columns.Item(6).ColumnWidth = 10

# Update the active selection / view
$ws.Range("C12").Select()
